# Threat-model "tag matrix" update: add DevOps-pipeline elements/technologies.
# The sheet is rebuilt in place: existing header/row-label/data formatting is
# extended (via copy/paste-special of formats only, so no new cell styles are
# created) to the new 12-column x 30-row extent, then every cell's value is
# (re)written to match the target layout. This naturally adds the new shared
# strings and reuses existing ones (e.g. "X") where the diff does.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Extend formatting to the new cells *before* touching values, by copying
#    formats (not values) from the existing reference cells/ranges:
#      - header row (B1 style) -> new header cells C1:L1
#      - row-label column (A2 style) -> new row labels A3:A30
#      - data-cell style (B2 style) -> the whole B:L data block incl. new rows/cols
#    Re-applying to cells that already carry that style is a no-op.
# ---------------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("C1:L1").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("A3:A30").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B2:L30").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. New columns H:L get the same 35-wide custom column width as B:G.
#    Excel's COM ColumnWidth is in "characters"; the stored OOXML width is
#    ColumnWidth + 5/6, so set ColumnWidth = 35 - 5/6 to land on width="35".
# ---------------------------------------------------------------------------
$newColWidth = 35 - (5/6)
for ($c = 8; $c -le 12; $c++) {
  $ws.Columns.Item($c).ColumnWidth = $newColWidth
}

# ---------------------------------------------------------------------------
# 3. Write every cell value for the final A1:L30 grid (row-major, A..L).
# ---------------------------------------------------------------------------
$data = @(
    @("Element","amazon ecr","amazon eks","aws secret manager","github","github action","mysql","nexus","nginx","react","spring","tomcat"),
    @("AWS Secret Manager Vault","","","X","","","","","","","",""),
    @("Amazon ECR Container Registry","X","","","","","","","","","",""),
    @("Amazon EKS Container Platform","","X","","","","","","","","",""),
    @("Container Platform Pull","","","","","","","","","","",""),
    @("Backend","","","","","","","","","","X","X"),
    @("Vault Access (backend)","","","","","","","","","","",""),
    @("Server Traffic","","","","","","","","","","",""),
    @("Database","","","","","","X","","","","",""),
    @("Development Client","","","","","","","","","","",""),
    @("Sourcecode Repository Traffic","","","","","","","","","","",""),
    @("Container Registry Traffic","","","","","","","","","","",""),
    @("Container Platform Traffic","","","","","","","","","","",""),
    @("Build Pipeline Traffic","","","","","","","","","","",""),
    @("Artifact Registry Traffic","","","","","","","","","","",""),
    @("Frontend","","","","","","","","X","X","",""),
    @("User Traffic","","","","","","","","","","",""),
    @("Nexus Artifact Registry","","","","","","","X","","","",""),
    @("github Sourcecode Repository","","","","X","","","","","","",""),
    @("github action Build Pipeline","","","","","X","","","","","",""),
    @("Sourcecode Repository Traffic","","","","","","","","","","",""),
    @("Container Registry Traffic","","","","","","","","","","",""),
    @("Artifact Registry Traffic","","","","","","","","","","",""),
    @("Configuration Secrets","","","","","","","","","","",""),
    @("Deployment","","","","","","","","","","",""),
    @("Greetings","","","","","","","","","","",""),
    @("Sourcecode","","","","","","","","","","",""),
    @("Trust Boundary","","","","","","","","","","",""),
    @("Amazon EKS Runtime","","X","","","","","","","","",""),
    @("EKS","","","","","","","","","","","")
)

for ($r = 0; $r -lt $data.Length; $r++) {
  $row = $data[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
  }
}
